# Update the practice-problem answers ("##÷#=") in the single table.
# Several old values are duplicated (e.g. "38÷2=" occurs three times with
# three different replacements), so we address each problem cell
# positionally via Table.Cell(row, col) rather than a global Find/Replace.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "27÷3="
$t.Cell(1, 2).Range.Text = "90÷2="
$t.Cell(1, 3).Range.Text = "18÷3="
$t.Cell(1, 4).Range.Text = "14÷8="
$t.Cell(1, 5).Range.Text = "48÷6="

# Row 5
$t.Cell(5, 1).Range.Text = "13÷5="
$t.Cell(5, 2).Range.Text = "40÷3="
$t.Cell(5, 3).Range.Text = "20÷9="
$t.Cell(5, 4).Range.Text = "97÷5="
$t.Cell(5, 5).Range.Text = "16÷8="

# Row 9
$t.Cell(9, 1).Range.Text = "71÷2="
$t.Cell(9, 2).Range.Text = "39÷6="
$t.Cell(9, 3).Range.Text = "54÷9="
$t.Cell(9, 4).Range.Text = "49÷8="
$t.Cell(9, 5).Range.Text = "71÷2="

# Row 13
$t.Cell(13, 1).Range.Text = "78÷9="
$t.Cell(13, 2).Range.Text = "49÷2="
$t.Cell(13, 3).Range.Text = "99÷7="
$t.Cell(13, 4).Range.Text = "91÷9="
$t.Cell(13, 5).Range.Text = "93÷7="

# Row 17
$t.Cell(17, 1).Range.Text = "32÷5="
$t.Cell(17, 2).Range.Text = "35÷3="
$t.Cell(17, 3).Range.Text = "35÷5="
$t.Cell(17, 4).Range.Text = "21÷6="
$t.Cell(17, 5).Range.Text = "41÷6="

Write-Output "Updated 25 problem cells."
